$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.608.99"
$ws.Range("D3").Value = "1.694.58"
$ws.Range("E3").Value = "  -5.70%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "219.78"
$ws.Range("E5").Value = "  -5.21%  "
$ws.Range("D6").Value = "0.5098"
$ws.Range("E6").Value = "  -13.29%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "0.2648"
$ws.Range("E8").Value = "  -4.36%  "
$ws.Range("D9").Value = "22.15"
$ws.Range("E9").Value = "  -4.71%  "
$ws.Range("D10").Value = "0.06291"
$ws.Range("E10").Value = "  -7.24%  "
$ws.Range("D11").Value = "0.07369"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").Value = "1.698.21"
$ws.Range("E12").Value = "  -5.35%  "
$ws.Range("D13").Value = "4.520"
$ws.Range("E13").Value = "  -5.92%  "
$ws.Range("D14").Value = "0.5794"
$ws.Range("E14").Value = "  -6.39%  "
$ws.Range("D15").Value = "1.925.90"
$ws.Range("E15").Value = "  -5.67%  "
$ws.Range("D16").Value = "0.000008449"
$ws.Range("E16").Value = "  -7.09%  "
$ws.Range("D17").Value = "65.49"
$ws.Range("E17").Value = "  -13.15%  "
$ws.Range("D18").Value = "26.634.29"
$ws.Range("E18").Value = "  -7.03%  "
$ws.Range("D19").Value = "4.987"
$ws.Range("E19").Value = "  -8.96%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "10.99"
$ws.Range("E21").Value = "  -4.58%  "
$ws.Range("D22").Value = "186.29"
$ws.Range("E22").Value = "  -11.53%  "
$ws.Range("D23").Value = "6.249"
$ws.Range("E23").Value = "  -8.56%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "144.68"
$ws.Range("E25").Value = "  -5.90%  "
$ws.Range("D26").Value = "7.498"
$ws.Range("E26").Value = "  -6.15%  "
$ws.Range("D27").Value = "0.1158"
$ws.Range("E27").Value = "  -8.31%  "
$ws.Range("D28").Value = "15.83"
$ws.Range("E28").Value = "  -3.78%  "
$ws.Range("D29").Value = "1.343"
$ws.Range("E29").Value = "  -5.52%  "
$ws.Range("D30").Value = "0.05684"
$ws.Range("E30").Value = "  -6.93%  "
$ws.Range("E31").Value = "  -6.34%  "
$ws.Range("D32").Value = "3.508"
$ws.Range("E32").Value = "  -7.41%  "
$ws.Range("D33").Value = "3.496"
$ws.Range("E33").Value = "  -8.09%  "
$ws.Range("D34").Value = "1.653"
$ws.Range("E34").Value = "  -4.97%  "
$ws.Range("D35").Value = "1.020"
$ws.Range("E35").Value = "  -2.82%  "
$ws.Range("D36").Value = "0.6002"
$ws.Range("E36").Value = "  -6.64%  "
$ws.Range("D37").Value = "2.360"
$ws.Range("E37").Value = "  -5.71%  "
$ws.Range("D38").Value = "2.692"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("D39").Value = "0.01618"
$ws.Range("E39").Value = "  -4.84%  "
$ws.Range("D40").Value = "1.101.85"
$ws.Range("E40").Value = "  -3.57%  "
$ws.Range("D41").Value = "0.8564"
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("D42").Value = "5.836"
$ws.Range("E42").Value = "  -10.27%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "99.41"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").Value = "1.851.75"
$ws.Range("E45").Value = "  -5.14%  "
$ws.Range("D46").Value = "0.00000000115"
$ws.Range("E46").Value = "  +2.27%  "
$ws.Range("D47").Value = "56.58"
$ws.Range("E47").Value = "  -6.04%  "
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").Value = "8.119"
$ws.Range("E49").Value = "  -2.82%  "
$ws.Range("E50").Value = "  -4.27%  "
$ws.Range("D51").Value = "0.4322"
$ws.Range("E51").Value = "  -3.54%  "
